# WasteReport update: replace the "Регистрационный номер призводителя"
# column header with "УНП", and move the active selection to E2:E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# Header cell D2 (merged D2:D3) changes text.
$ws.Range("D2").Value = "УНП"

# Update the active selection shown in the saved view.
$ws.Activate()
$ws.Range("E2:E3").Select()
